$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo: "Phased Array (Deploable)" -> "Phased Array (Deployable)"
$ws.Range("A3").Value = "Phased Array (Deployable)"

# Add new rows for Multi Bandwidth Dish Transceiver variants
$ws.Range("A10").Value = "Multi Bandwidth Dish Transceiver (large)"
$ws.Range("C10").Value = 20
$ws.Range("F10").Value = "yes"
$ws.Range("G10").Value = 24
$ws.Range("H10").Value = 21000

$ws.Range("A11").Value = "Multi Bandwidth Dish Transceiver (medium)"
$ws.Range("C11").Value = 10
$ws.Range("F11").Value = "yes"
$ws.Range("G11").Value = 8
$ws.Range("H11").Value = 7000

$ws.Range("A12").Value = "Multi Bandwidth Dish Transceiver (shielded)"
$ws.Range("C12").Value = 5
$ws.Range("F12").Value = "yes"
$ws.Range("G12").Value = 4
$ws.Range("H12").Value = 3500

# Match the author's final cursor/selection position
$ws.Range("L6").Select()
